$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("median_age", "estimate_total_total_population_summary_indicators_median_age_years"),
    @("median_age_male", "estimate_male_total_population_summary_indicators_median_age_years"),
    @("median_age_female", "estimate_female_total_population_summary_indicators_median_age_years"),
    @("sex_ratio", "estimate_total_total_population_summary_indicators_sex_ratio_males_per_100_females"),
    @("age_dependency", "estimate_total_total_population_summary_indicators_age_dependency_ratio"),
    @("old_age_dependency", "estimate_total_total_population_summary_indicators_old_age_dependency_ratio"),
    @("child_dependency", "estimate_total_total_population_summary_indicators_child_dependency_ratio")
)

$startRow = 115
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}

# Set column A width
$ws.Columns.Item(1).ColumnWidth = 20.33203125

# Update the view: top-left cell, zoom, and selection
$ws.Application.ActiveWindow.ScrollRow = 102
$ws.Application.ActiveWindow.Zoom = 194
$ws.Range("B119").Select()
